$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 630, pushing the existing rows 630:671 down to 631:672.
$ws.Rows.Item(630).Insert()

# Column A holds a date-formatted string ("yyyy/mm/dd"); force it to stay
# plain text instead of Excel auto-converting it to a date serial number.
$ws.Range("A630").NumberFormat = "@"
$ws.Range("A630").Value = "2026/01/15"
$ws.Range("A630").Style = "Normal"

$ws.Range("B630").Value = "木"
$ws.Range("C630").Value = 16
$ws.Range("D630").Value = 35
